# Updates odds/score data cells on Sheet1 to match the 2025-04-15 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("N2").Value = 4.2
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 1.95
$ws.Range("Q2").Value = 1.85

# Row 3
$ws.Range("G3").Value = 2.1
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 3.9
$ws.Range("J3").Value = 1.1
$ws.Range("K3").Value = 7
$ws.Range("N3").Value = 2.4
$ws.Range("O3").Value = 1.53
$ws.Range("U3").Value = 9
$ws.Range("Z3").Value = 7
$ws.Range("AA3").Value = 6
$ws.Range("AB3").Value = 17
$ws.Range("AE3").Value = 9

# Row 4
$ws.Range("G4").Value = 2.8
$ws.Range("I4").Value = 2.63
$ws.Range("J4").Value = 1.1
$ws.Range("K4").Value = 7
$ws.Range("N4").Value = 2.4
$ws.Range("O4").Value = 1.53
$ws.Range("W4").Value = 29
$ws.Range("X4").Value = 26
$ws.Range("AE4").Value = 7
$ws.Range("AH4").Value = 26

# Row 5
$ws.Range("G5").Value = 3.8
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 1.8
$ws.Range("R5").Value = 1.53
$ws.Range("S5").Value = 2.38
$ws.Range("T5").Value = 15
$ws.Range("U5").Value = 23
$ws.Range("V5").Value = 13
$ws.Range("W5").Value = 41
$ws.Range("X5").Value = 29
$ws.Range("Y5").Value = 29
$ws.Range("AA5").Value = 8
$ws.Range("AB5").Value = 13
$ws.Range("AC5").Value = 41
$ws.Range("AD5").Value = 126
$ws.Range("AE5").Value = 10
$ws.Range("AF5").Value = 11
$ws.Range("AG5").Value = 8.5
$ws.Range("AH5").Value = 15
$ws.Range("AI5").Value = 13
$ws.Range("AJ5").Value = 19

# Row 6
$ws.Range("G6").Value = 4.75
$ws.Range("H6").Value = 3.9
$ws.Range("L6").Value = 1.22
$ws.Range("M6").Value = 4
$ws.Range("N6").Value = 1.7
$ws.Range("O6").Value = 2.1
$ws.Range("P6").Value = 1.33
$ws.Range("Q6").Value = 3.25
$ws.Range("R6").Value = 1.7
$ws.Range("S6").Value = 2.05
$ws.Range("X6").Value = 34
$ws.Range("Y6").Value = 34
$ws.Range("Z6").Value = 13
$ws.Range("AA6").Value = 7.5
$ws.Range("AC6").Value = 41
$ws.Range("AD6").Value = 151
$ws.Range("AE6").Value = 8.5
$ws.Range("AF6").Value = 9

# Row 7
$ws.Range("N7").Value = 1.6
$ws.Range("O7").Value = 2.3
$ws.Range("T7").Value = 9.5
$ws.Range("U7").Value = 10
$ws.Range("Z7").Value = 17
$ws.Range("AG7").Value = 13
$ws.Range("AH7").Value = 41
$ws.Range("AI7").Value = 29
$ws.Range("AJ7").Value = 29

# Row 8
$ws.Range("G8").Value = 2.15
$ws.Range("H8").Value = 2.75
$ws.Range("I8").Value = 4.33
$ws.Range("J8").Value = 1.14
$ws.Range("K8").Value = 5.5
$ws.Range("N8").Value = 3.1
$ws.Range("O8").Value = 1.36
$ws.Range("P8").Value = 1.73
$ws.Range("Q8").Value = 2
$ws.Range("U8").Value = 8.5
$ws.Range("AE8").Value = 8
$ws.Range("AF8").Value = 19
$ws.Range("AJ8").Value = 67

# Row 9
$ws.Range("G9").Value = 2.6
$ws.Range("H9").Value = 2.8
$ws.Range("I9").Value = 3.1
$ws.Range("J9").Value = 1.14
$ws.Range("K9").Value = 5.5
$ws.Range("L9").Value = 1.57
$ws.Range("M9").Value = 2.25
$ws.Range("N9").Value = 2.88
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 1.62
$ws.Range("Q9").Value = 2.2
$ws.Range("R9").Value = 2.25
$ws.Range("S9").Value = 1.57
$ws.Range("U9").Value = 11
$ws.Range("V9").Value = 11
$ws.Range("W9").Value = 26
$ws.Range("X9").Value = 26
$ws.Range("Z9").Value = 5.5
$ws.Range("AA9").Value = 5.5
$ws.Range("AC9").Value = 81
$ws.Range("AE9").Value = 7
$ws.Range("AF9").Value = 13

# Row 13
$ws.Range("K13").Value = 8
$ws.Range("N13").Value = 2.25
$ws.Range("O13").Value = 1.62

# Row 14
$ws.Range("N14").Value = 2.03
$ws.Range("O14").Value = 1.78

# Row 15
$ws.Range("H15").Value = 3.25
$ws.Range("I15").Value = 2.35
$ws.Range("N15").Value = 2.05
$ws.Range("O15").Value = 1.75
$ws.Range("P15").Value = 1.44
$ws.Range("Q15").Value = 2.63
$ws.Range("R15").Value = 1.8
$ws.Range("S15").Value = 1.91
$ws.Range("T15").Value = 9
$ws.Range("X15").Value = 26
$ws.Range("Z15").Value = 9
$ws.Range("AA15").Value = 6
$ws.Range("AE15").Value = 7.5
$ws.Range("AH15").Value = 23
$ws.Range("AI15").Value = 21

# Row 16
$ws.Range("H16").Value = 6
$ws.Range("S16").Value = 1.46
$ws.Range("Z16").Value = 15
$ws.Range("AB16").Value = 28
$ws.Range("AC16").Value = 120
$ws.Range("AE16").Value = 6.7
$ws.Range("AF16").Value = 5
$ws.Range("AG16").Value = 8.75
$ws.Range("AH16").Value = 5.3
$ws.Range("AI16").Value = 9.25

# Row 17
$ws.Range("H17").Value = 3.3
$ws.Range("I17").Value = 2.92
$ws.Range("L17").Value = 1.26
$ws.Range("M17").Value = 3.45
$ws.Range("N17").Value = 1.83
$ws.Range("O17").Value = 1.78
$ws.Range("P17").Value = 1.37
$ws.Range("Q17").Value = 2.5
$ws.Range("R17").Value = 1.69
$ws.Range("S17").Value = 2.04
$ws.Range("T17").Value = 6.9
$ws.Range("U17").Value = 9
$ws.Range("W17").Value = 17
$ws.Range("X17").Value = 14
$ws.Range("Y17").Value = 21
$ws.Range("Z17").Value = 10
$ws.Range("AA17").Value = 5.6
$ws.Range("AB17").Value = 11.25
$ws.Range("AC17").Value = 45
$ws.Range("AD17").Value = 300
$ws.Range("AE17").Value = 7.9
$ws.Range("AF17").Value = 12.5
$ws.Range("AG17").Value = 9
$ws.Range("AI17").Value = 20
$ws.Range("AJ17").Value = 26

# Row 18
$ws.Range("G18").Value = 1.9
$ws.Range("H18").Value = 3.3
$ws.Range("I18").Value = 4.33
$ws.Range("K18").Value = 9
$ws.Range("R18").Value = 1.95
$ws.Range("S18").Value = 1.8
$ws.Range("T18").Value = 6.5
$ws.Range("U18").Value = 8.5
$ws.Range("V18").Value = 8.5
$ws.Range("W18").Value = 15
$ws.Range("X18").Value = 17
$ws.Range("AA18").Value = 6.5
$ws.Range("AE18").Value = 11
$ws.Range("AF18").Value = 21
$ws.Range("AG18").Value = 15
$ws.Range("AI18").Value = 34

# Row 19
$ws.Range("G19").Value = 1.42
$ws.Range("H19").Value = 4.75
$ws.Range("L19").Value = 1.22
$ws.Range("M19").Value = 4
$ws.Range("N19").Value = 1.75
$ws.Range("O19").Value = 2.05
$ws.Range("P19").Value = 1.33
$ws.Range("Q19").Value = 3.25
$ws.Range("R19").Value = 2
$ws.Range("S19").Value = 1.75
$ws.Range("T19").Value = 7
$ws.Range("X19").Value = 12
$ws.Range("Y19").Value = 26
$ws.Range("Z19").Value = 13
$ws.Range("AA19").Value = 9
$ws.Range("AD19").Value = 351
$ws.Range("AE19").Value = 17

# Row 20
$ws.Range("N20").Value = 1.88
$ws.Range("O20").Value = 1.93
